# Add data for 2022-06-19
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update header label to reflect the new "through" date
$ws.Name = "Through 2022-06-11"
$ws.Range("B1").Value = "June 2022 (through June 11)"

# Update existing counts
$ws.Range("T2").Value = 2
$ws.Range("T3").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("B10").Value = 4
$ws.Range("AL10").Value = 3
$ws.Range("N14").Value = 2
$ws.Range("N38").Value = 2
$ws.Range("B92").Value = 2

# Add new counts to previously-empty cells
$ws.Range("N12").Value = 1
$ws.Range("AR14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("Z19").Value = 1
$ws.Range("B27").Value = 1
$ws.Range("AL27").Value = 1
$ws.Range("AF37").Value = 1
$ws.Range("N68").Value = 1
$ws.Range("H70").Value = 1
$ws.Range("B71").Value = 1
